$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.052.91'
$ws.Range("E2").Value = '  +3.12%  '
$ws.Range("D3").Value = '3.803.94'
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''699.67'
$ws.Range("E5").Value = '  +8.01%  '
$ws.Range("D6").Value = '''172.93'
$ws.Range("E6").Value = '  +4.58%  '
$ws.Range("D7").Value = '3.802.05'
$ws.Range("E7").Value = '  +1.10%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +1.25%  '
$ws.Range("E11").Value = '  +6.74%  '
$ws.Range("E12").Value = '  +1.29%  '
$ws.Range("E13").Value = '  +8.37%  '
$ws.Range("D14").Value = '''36.43'
$ws.Range("E14").Value = '  +4.68%  '
$ws.Range("D15").Value = '4.446.77'
$ws.Range("E15").Value = '  +1.22%  '
$ws.Range("D16").Value = '3.827.00'
$ws.Range("E16").Value = '  +1.44%  '
$ws.Range("D17").Value = '70.982.31'
$ws.Range("E17").Value = '  +3.09%  '
$ws.Range("D18").Value = '''17.91'
$ws.Range("E18").Value = '  +1.56%  '
$ws.Range("E19").Value = '  +3.03%  '
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("D22").Value = '''481.79'
$ws.Range("E22").Value = '  +3.21%  '
$ws.Range("E23").Value = '  +1.54%  '
$ws.Range("D24").Value = '''84.28'
$ws.Range("E24").Value = '  +3.04%  '
$ws.Range("E25").Value = '  +0.58%  '
$ws.Range("D26").Value = '''12.45'
$ws.Range("E26").Value = '  +2.17%  '
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '''10.54'
$ws.Range("E27").Value = '  +3.38%  '
$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").Value = '''2.18'
$ws.Range("E28").Value = '  +3.77%  '
$ws.Range("D29").Value = '3.955.97'
$ws.Range("E29").Value = '  +1.21%  '
$ws.Range("E30").Value = '  -0.09%  '
$ws.Range("E31").Value = '  +15.31%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '''2.30'
$ws.Range("E32").Value = '  +2.05%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").Value = '''7.56'
$ws.Range("E33").Value = '  +6.37%  '
$ws.Range("D34").Value = '''0.188'
$ws.Range("E34").Value = '  +9.26%  '
$ws.Range("D35").Value = '''29.51'
$ws.Range("E35").Value = '  +3.70%  '
$ws.Range("D36").Value = '''9.25'
$ws.Range("E36").Value = '  +5.29%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("E38").Value = '  +2.86%  '
$ws.Range("E39").Value = '  +7.53%  '
$ws.Range("D40").Value = '''6.03'
$ws.Range("E40").Value = '  +4.67%  '
$ws.Range("E41").Value = '  +12.93%  '
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").Value = '''0.977'
$ws.Range("E42").Value = '  +2.14%  '
$ws.Range("B43").Value = 'FLOKI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D43").Value = '''0.000327'
$ws.Range("E43").Value = '  +23.62%  '
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D46").Value = '''162.52'
$ws.Range("E46").Value = '  +4.89%  '
$ws.Range("D47").Value = '''44.85'
$ws.Range("E47").Value = '  -0.30%  '
$ws.Range("D48").Value = '''48.74'
$ws.Range("E48").Value = '  +3.12%  '
$ws.Range("E49").Value = '  +2.74%  '
$ws.Range("E50").Value = '  -0.59%  '
$ws.Range("D51").Value = '''8.58'
$ws.Range("E51").Value = '  +2.85%  '
